$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 35, shifting existing rows 35-139 down to 36-140.
$ws.Rows.Item(35).EntireRow.Insert()

# Populate the newly inserted row 35 with the new weekly data point.
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35, 3).Value = "Ñuble"
$ws.Cells.Item(35, 4).Value = 44498
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(35, 6).Value = 100112006
$ws.Cells.Item(35, 7).Value = "Repollo"
$ws.Cells.Item(35, 8).Value = "Crespo record"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 300
$ws.Cells.Item(35, 11).Value = 650
$ws.Cells.Item(35, 12).Value = 700
$ws.Cells.Item(35, 13).Value = 675
$ws.Cells.Item(35, 14).Value = "$/unidad"
$ws.Cells.Item(35, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(35, 16).Value = 675
$ws.Cells.Item(35, 17).Value = 1
$ws.Cells.Item(35, 18).Value = "Hortaliza"
